{"js": "// Fix a typo (\"Rewire\" -> \"Rewrite\") in the section heading, and change\n// \"years\" to \"days\" in the exponential-decay word problem.\n\n// 1) \"Rewire as from Exponential...\" -> \"Rewrite as from Exponential...\"\nconst headingResults = context.document.body.search(\"Rewire\", { matchCase: true });\nheadingResults.load(\"text\");\nawait context.sync();\n\nif (headingResults.items.length > 0) {\n  headingResults.items[0].insertText(\"Rewrite\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"...t is the time in years find the following...\" -> \"...in days find...\"\nconst yearsResults = context.document.body.search(\"years\", { matchCase: true });\nyearsResults.load(\"text\");\nawait context.sync();\n\nif (yearsResults.items.length > 0) {\n  yearsResults.items[0].insertText(\"days\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Fix a typo (\"Rewire\" -> \"Rewrite\") in the section heading, and change\n# \"years\" to \"days\" in the exponential-decay word problem.\n\n$d = $word.ActiveDocument\n\n# 1) \"Rewire as from Exponential...\" -> \"Rewrite as from Exponential...\"\n$find1 = $d.Content.Find\n$find1.Text = \"Rewire\"\n$find1.Replacement.Text = \"Rewrite\"\n$find1.Execute(\n  $find1.Text,              # FindText\n  $false,                   # MatchCase\n  $false,                   # MatchWholeWord\n  $false,                   # MatchWildcards\n  $false,                   # MatchSoundsLike\n  $false,                   # MatchAllWordForms\n  $true,                    # Forward\n  1,                        # Wrap (wdFindContinue)\n  $false,                   # Format\n  $find1.Replacement.Text,  # ReplaceWith\n  2                         # Replace (wdReplaceAll)\n)\n\n# 2) \"...t is the time in years find the following...\" -> \"...in days find...\"\n$find2 = $d.Content.Find\n$find2.Text = \"years\"\n$find2.Replacement.Text = \"days\"\n$find2.Execute(\n  $find2.Text,              # FindText\n  $false,                   # MatchCase\n  $false,                   # MatchWholeWord\n  $false,                   # MatchWildcards\n  $false,                   # MatchSoundsLike\n  $false,                   # MatchAllWordForms\n  $true,                    # Forward\n  1,                        # Wrap (wdFindContinue)\n  $false,                   # Format\n  $find2.Replacement.Text,  # ReplaceWith\n  2                         # Replace (wdReplaceAll)\n)\n"}
